$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.4761386265337256
$ws.Range("D2").Value = 0.6386701084712101

$ws.Range("C3").Value = 0.2636000925935691
$ws.Range("D3").Value = 0.7945413848735541

$ws.Range("C4").Value = 0.2890216310624869
$ws.Range("D4").Value = 0.7752707940813921

$ws.Range("C5").Value = 0.8130170337770406
$ws.Range("D5").Value = 0.4249198429656413

$ws.Range("C6").Value = -0.1270425824062188
$ws.Range("D6").Value = 0.9000607224754804

$ws.Range("C7").Value = -0.1445293758247178
$ws.Range("D7").Value = 0.8863983175187313

$ws.Range("C8").Value = 0.2090958313990981
$ws.Range("D8").Value = 0.8362974275745216

$ws.Range("C9").Value = 0.02176825797007874
$ws.Range("D9").Value = 0.9828290465389999

$ws.Range("C10").Value = 0.2452404788745877
$ws.Range("D10").Value = 0.8085434160545089

$ws.Range("C11").Value = 0.2795455790027953
$ws.Range("D11").Value = 0.7824374935313334
